$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header cell H1 "Save", matching the formatting of the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the Save column values for rows 2-4
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
